$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 157, shifting existing rows 157:168 down to 158:169
$ws.Rows("157:157").Insert()

# Populate the newly inserted row 157 with the new weekly data record
$ws.Range("A157").Value = 10
$ws.Range("B157").Value = "Vega Modelo de Temuco"
$ws.Range("C157").Value = "La Araucanía"
$ws.Range("D157").Value = 44578
$ws.Range("E157").Value = 9
$ws.Range("F157").Value = 100112043
$ws.Range("G157").Value = "Pepino dulce"
$ws.Range("H157").Value = "Cultivar IV Región"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 300
$ws.Range("K157").Value = 24000
$ws.Range("L157").Value = 24000
$ws.Range("M157").Value = 24000
$ws.Range("N157").Value = '$/bandeja 18 kilos'
$ws.Range("O157").Value = "Provincia de Limarí"
$ws.Range("P157").Value = 1333
$ws.Range("Q157").Value = 18
$ws.Range("R157").Value = "Hortaliza"

# Ensure the date cell keeps the same style/number format as the other date cells in column D
$ws.Range("D157").NumberFormat = $ws.Range("D158").NumberFormat
